# Update countries & provincias Spain
# - refresh the "Datos actualizados..." timestamp banner
# - refresh case counters for Estados Unidos, Brasil, Alemania, Uruguay

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp banner (row 1)
$ws.Range("A1").Value = "Datos actualizados a 18 de Mayo de 2020 a las 00:35"

# Estados Unidos (row 4): Casos totales, Nuevos casos, Casos activos, Recuperados, Muertes hoy, Muertes
$ws.Range("B4").Value = 1526136
$ws.Range("C4").Value = 18363
$ws.Range("D4").Value = 344805
$ws.Range("E4").Value = 1090400
$ws.Range("G4").Value = 818
$ws.Range("H4").Value = 90931

# Brasil (row 8)
$ws.Range("B8").Value = 240307
$ws.Range("C8").Value = 7165
$ws.Range("E8").Value = 134544
$ws.Range("G8").Value = 458
$ws.Range("H8").Value = 16091

# Alemania (row 11)
$ws.Range("B11").Value = 176651
$ws.Range("C11").Value = 407
$ws.Range("E11").Value = 15202
$ws.Range("G11").Value = 22
$ws.Range("H11").Value = 8049

# Uruguay (row 118)
$ws.Range("B118").Value = 734
$ws.Range("C118").Value = 1
$ws.Range("D118").Value = 564
$ws.Range("E118").Value = 150
$ws.Range("G118").Value = 1
$ws.Range("H118").Value = 20
